# The crawl re-ran, refreshing the "timestamp" column (O) for every
# scraped row (rows 2-33) on Sheet1 with the new crawl time.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2:O33").Value = "2022-08-29 20:59:05"
